{"js": "// Wrap the final \"project description\" sentence in a C-style comment\n// block and extend it with the extra sentences about finishing the\n// initiation phase (background description, technical parts, risk\n// assessment, layout, checking and presenting).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetText = \"The last part was writing the project description.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetText) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'The last part was writing...' paragraph\");\n}\n\n// New paragraph right before the target, opening the comment block.\ntarget.insertParagraph(\n  \"/* Maybe in a totally different way, it\\u2019s been to many hours to think about it now\\u2026\",\n  \"Before\"\n);\n\n// New paragraph right after the target, closing the comment block.\ntarget.insertParagraph(\"*/\", \"After\");\n\n// Append the extra sentences to the original sentence, inside the same\n// paragraph (kept before the _GoBack bookmark, at the paragraph's end).\ntarget.insertText(\" \", \"End\");\ntarget.insertText(\n  \"We gathered the ideas for background description together and then Daniela put them into words. Michael took care of the more technical parts and Matej created the risk assessment section and gathered everything in the correct layout. The final part was checking it by everyone. \",\n  \"End\"\n);\ntarget.insertText(\n  \"Having it done, we had to decide how to present it. Michael with Michaela took presenting and making the presentation and Matej with Daniela took care of the feedback.\",\n  \"End\"\n);\n\nawait context.sync();\n", "ps1": "# Wrap the final \"project description\" sentence in a C-style comment\n# block and extend it with the extra sentences about finishing the\n# initiation phase (background description, technical parts, risk\n# assessment, layout, checking and presenting).\n\n$d = $word.ActiveDocument\n\n$targetText = \"The last part was writing the project description.\"\n\nfunction Find-TargetIndex {\n    $i = 0\n    foreach ($p in $d.Paragraphs) {\n        $i = $i + 1\n        if ($p.Range.Text.TrimEnd(\"`r\") -eq $targetText) {\n            return $i\n        }\n    }\n    return -1\n}\n\n$idx = Find-TargetIndex\nif ($idx -eq -1) {\n    throw \"Could not find the 'The last part was writing...' paragraph\"\n}\n\n# Insert a new (empty) paragraph right before the target and fill it with\n# the comment-opening line.\n$target = $d.Paragraphs.Item($idx)\n$target.Range.InsertParagraphBefore()\n\n$idx = Find-TargetIndex\n$openPara = $d.Paragraphs.Item($idx - 1)\n$openPara.Range.Text = \"/* Maybe in a totally different way, it\u2019s been to many hours to think about it now\u2026\"\n\n# Insert a new (empty) paragraph right after the target and fill it with\n# the comment-closing marker.\n$target = $d.Paragraphs.Item($idx)\n$target.Range.InsertParagraphAfter()\n$closePara = $d.Paragraphs.Item($idx + 1)\n$closePara.Range.Text = \"*/\"\n\n# Append the extra sentences to the original sentence, inside the same\n# paragraph (kept before the _GoBack bookmark, at the paragraph's end).\n$target = $d.Paragraphs.Item($idx)\n$target.Range.InsertAfter(\" \")\n$target.Range.InsertAfter(\"We gathered the ideas for background description together and then Daniela put them into words. Michael took care of the more technical parts and Matej created the risk assessment section and gathered everything in the correct layout. The final part was checking it by everyone. \")\n$target.Range.InsertAfter(\"Having it done, we had to decide how to present it. Michael with Michaela took presenting and making the presentation and Matej with Daniela took care of the feedback.\")\n"}
